# JSExecutor Use cases Incorporated
#
# US2-Delete gets two more rows appended, reusing the already-known
# "Updated ..." filter names (existing shared strings), formatted like the
# other data rows (thin border, centered) but without vertical centering.
# Selection/active-tab bookkeeping moves from US3-Edit to US2-Delete.

$wb = $excel.ActiveWorkbook

# --- US2-Delete: append the two new "deleted/updated filter" rows ---
$ws2 = $wb.Worksheets.Item("US2-Delete")
$ws2.Range("A4").Value = "Updated Auto Bugs Filter"
$ws2.Range("A5").Value = "Updated Android Filter"

# Match the look of the existing data rows (border + centered text), then
# drop the vertical centering the source row has.
$fmtSource = $ws2.Range("A2")
[void]$fmtSource.Copy()
$newRows = $ws2.Range("A4:A5")
[void]$newRows.PasteSpecial(-4122)
$newRows.VerticalAlignment = -4107

# --- US3-Edit: selection moves to B3, tab no longer active ---
$ws3 = $wb.Worksheets.Item("US3-Edit")
[void]$ws3.Activate()
[void]$ws3.Range("B3").Select()

# --- US2-Delete becomes the active tab/sheet, selection moves to E5 ---
[void]$ws2.Activate()
[void]$ws2.Range("E5").Select()
